$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.465899999999994
$ws.Range("D21").Value = -8.395499999999997
$ws.Range("D23").Value = -8.083699999999997
$ws.Range("D25").Value = -7.346799999999995
$ws.Range("D53").Value = -8.076899999999997
$ws.Range("D57").Value = -8.298699999999995
$ws.Range("D59").Value = -8.221099999999993
$ws.Range("D69").Value = -7.331499999999997
$ws.Range("D79").Value = -7.306200000000009
$ws.Range("D83").Value = -9.164699999999996
$ws.Range("D93").Value = -7.068199999999993
